$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<then>"
$ws.Range("C2").Value = 48

$ws.Range("B4").Value = "<sierra>"
$ws.Range("C4").Value = 57

$ws.Range("C5").Value = 53

$ws.Range("B6").Value = "<part>"
$ws.Range("C6").Value = 51

$ws.Range("C7").Value = 47

$ws.Range("C8").Value = 49

$ws.Range("B9").Value = "<when>"
$ws.Range("C9").Value = 48

$ws.Range("B10").Value = "<paste>"

$ws.Range("C11").Value = 48

$ws.Range("B12").Value = "<into>"
$ws.Range("C12").Value = 55

$ws.Range("B13").Value = "<uniform>"
$ws.Range("C13").Value = 59

$ws.Range("B14").Value = "<for>"
$ws.Range("C14").Value = 56

$ws.Range("B15").Value = "<more>"
$ws.Range("C15").Value = 56

$ws.Range("C16").Value = 27
